# Backup before MoClo code restructure
# Expand the destination-well pattern from A1..A4 (4 wells) to A1..A6 (6 wells)
# for each of the 3 source wells (A1, A2, A3), growing the data block from
# rows 2-13 (12 rows) to rows 2-19 (18 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B-F, and the Volume/Reagent pairing per source-well group, stay the
# same as before; only the number of destination wells per group changes
# from 4 to 6, and the row count grows accordingly.

$sourcePlateName = "level 1 LDV source plate"
$sourcePlateType = "384LDV_AQ_B"
$destPlateName    = "384-Well Level 1 MoClo output plate"
$destPlateType    = "Echo® Qualified 384-Well Polypropylene Source Microplate (384PP)"

$sourceWells = @("A1", "A2", "A3")
$destWells   = @("A1", "A2", "A3", "A4", "A5", "A6")
$volumes     = @{ "A1" = 500; "A2" = 125; "A3" = 250 }
$reagents    = @{ "A1" = "DNA ligase buffer"; "A2" = "DNA ligase"; "A3" = "BsaI-HF" }

$uid = 1
$row = 2

foreach ($sw in $sourceWells) {
    foreach ($dw in $destWells) {
        $ws.Cells.Item($row, 1).Value = $uid
        $ws.Cells.Item($row, 2).Value = $sourcePlateName
        $ws.Cells.Item($row, 3).Value = $sourcePlateType
        $ws.Cells.Item($row, 4).Value = $sw
        $ws.Cells.Item($row, 5).Value = $destPlateName
        $ws.Cells.Item($row, 6).Value = $destPlateType
        $ws.Cells.Item($row, 7).Value = $dw
        $ws.Cells.Item($row, 8).Value = $volumes[$sw]
        $ws.Cells.Item($row, 9).Value = $reagents[$sw]

        $uid = $uid + 1
        $row = $row + 1
    }
}
